$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.818.49"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "3.327.43"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'581.69"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'174.97"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").Value = "3.323.52"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").Value = "'0.182"
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("D11").Value = "'0.579"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "'47.02"
$ws.Range("E12").Value = "  +4.06%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").Value = "'696.65"
$ws.Range("E14").Value = "  +4.37%  "
$ws.Range("D15").Value = "3.861.68"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "'8.36"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").Value = "67.755.38"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "3.331.19"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "'11.09"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").Value = "'0.889"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Value = "'5.41"
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "'16.89"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'101.33"
$ws.Range("E25").Value = "  +2.99%  "
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "'2.68"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").Value = "'9.40"
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("D29").Value = "'32.93"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").Value = "'8.52"
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").Value = "'6.99"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").Value = "'572.26"
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").Value = "'10.98"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "3.717.30"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'56.53"
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("D38").Value = "'3.28"
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("D39").Value = "'35.48"
$ws.Range("E39").Value = "  +9.85%  "
$ws.Range("D40").Value = "'0.133"
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("D41").Value = "'3.13"
$ws.Range("E41").Value = "  +2.52%  "
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").Value = "'3.32"
$ws.Range("E43").Value = "  +2.82%  "
$ws.Range("D44").Value = "0.0₃0670"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").Value = "'0.334"
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").Value = "'0.0405"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("D51").Value = "'130.99"
$ws.Range("E51").Value = "  +1.59%  "
